$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status moves from "In Translation" to "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Latest HO Xliff Generate Date / Latest Handoff Datetime bump forward a few seconds
$wsOverview.Range("G2").Value = "2016-09-05 22:50:43"
$wsDeDe.Range("H2").Value = "2016-09-05 22:50:43"
$wsZhCn.Range("H2").Value = "2016-09-05 22:50:39"

# Widen the Status columns to fit "Ready for handoff"
$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797
$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797
